# Update teacher names in column A (name) for Sheet1
# Each row's name value is updated to match the new value from the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(36, 1).Value = 'سید مصطفی سیادت موسوی'
$ws.Cells.Item(52, 1).Value = 'سید محمد شهرتاش'
$ws.Cells.Item(54, 1).Value = 'سید محمدرضا موسوی میرکلایی'
$ws.Cells.Item(113, 1).Value = 'سید حسن هاشم آبادی'
$ws.Cells.Item(114, 1).Value = 'سید محمد'
$ws.Cells.Item(130, 1).Value = 'سیده اشم مسدد'
$ws.Cells.Item(133, 1).Value = 'سید سعید سیادت نژاد'
$ws.Cells.Item(160, 1).Value = 'سید اصغر موسوی'
$ws.Cells.Item(201, 1).Value = 'سید سینا ثمره موسوی'
$ws.Cells.Item(202, 1).Value = 'سید مهدی عطیفه کمال ابادفراهان'
$ws.Cells.Item(269, 1).Value = 'سید سپهر موسوی'
$ws.Cells.Item(270, 1).Value = 'سید علی نیک نام'
$ws.Cells.Item(283, 1).Value = 'سید مهدی علوی املشی'
$ws.Cells.Item(294, 1).Value = 'سید سجاد میرولد'
$ws.Cells.Item(350, 1).Value = 'سید جواد امامی'
$ws.Cells.Item(353, 1).Value = 'سید حسن صدیقی'
$ws.Cells.Item(386, 1).Value = 'سید صالح اعتمادی'
$ws.Cells.Item(394, 1).Value = 'سید مجید هاشمیان زاده'
$ws.Cells.Item(443, 1).Value = 'سید حجت سبزپوشان'
$ws.Cells.Item(459, 1).Value = 'سید حسین حسینی شکوه'
$ws.Cells.Item(473, 1).Value = 'سید محمدمهدی عبیری'
$ws.Cells.Item(488, 1).Value = 'سید میعاد صالحی'
$ws.Cells.Item(510, 1).Value = 'سید حسین'
$ws.Cells.Item(528, 1).Value = 'سید عبدالهادی دانشپور'
$ws.Cells.Item(603, 1).Value = 'سید اقایی'
$ws.Cells.Item(605, 1).Value = 'سید مصطفی حسینعلی پور'
$ws.Cells.Item(616, 1).Value = 'سید ادریس فیض آبادی'
$ws.Cells.Item(623, 1).Value = 'سید محمدعلی بوترابی'
$ws.Cells.Item(675, 1).Value = 'سید جواد ازهری'
$ws.Cells.Item(676, 1).Value = 'سید نظام الدین اشرفی زاده'
$ws.Cells.Item(681, 1).Value = 'سید محمدعلی موسوی'
$ws.Cells.Item(700, 1).Value = 'سید ذبیح الله طباطبایی شیرازانی'
$ws.Cells.Item(721, 1).Value = 'سید محمود میرطباطبایی'
$ws.Cells.Item(755, 1).Value = 'سید مجید مفیدی شمیرانی'
$ws.Cells.Item(837, 1).Value = 'سید مجتبی حسینی نسب'
$ws.Cells.Item(871, 1).Value = 'سید اصغر ابن الرسول'
$ws.Cells.Item(878, 1).Value = 'سید امیر منصوری'
$ws.Cells.Item(916, 1).Value = 'سید رضا علمی حسینی'
$ws.Cells.Item(917, 1).Value = 'سید مجتبی میرحسینی'
$ws.Cells.Item(932, 1).Value = 'سید علی حسینی'
$ws.Cells.Item(974, 1).Value = 'سید مهدی حسینی دولت آبادی'
$ws.Cells.Item(980, 1).Value = 'سید علی قهاری'
$ws.Cells.Item(985, 1).Value = 'سید حسن موسوی'
$ws.Cells.Item(988, 1).Value = 'سید حامد رستگار'
$ws.Cells.Item(1007, 1).Value = 'سید ابوالفضل حسینی زاده'
$ws.Cells.Item(1012, 1).Value = 'سید علیرضا'
$ws.Cells.Item(1030, 1).Value = 'سید دانیال غفاریان تربتی مجاور'
$ws.Cells.Item(1037, 1).Value = 'سید علی سرکشیکیان'
$ws.Cells.Item(1070, 1).Value = 'سید حمید حاجی'
$ws.Cells.Item(1080, 1).Value = 'سید حسن ذبیحی فر'
